# Scheduled runner refresh of market price / profit data across the
# Kraken_Profits sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Columns: H=currentAveragePrice, I=currentAveragePriceNQ,
# J=currentAveragePriceHQ, K=LevePriceNQ, L=LevePriceHQ,
# M=LeveProfitNQ, N=LeveProfitHQ.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(4, 8).Value = 2461
$ws.Cells.Item(4, 9).Value = 2461
$ws.Cells.Item(4, 11).Value = 2461
$ws.Cells.Item(4, 13).Value = -2347

$ws.Cells.Item(8, 8).Value = 1363.2
$ws.Cells.Item(8, 9).Value = 20
$ws.Cells.Item(8, 11).Value = 60
$ws.Cells.Item(8, 13).Value = 79

$ws.Cells.Item(19, 8).Value = 498.1111
$ws.Cells.Item(19, 9).Value = 497.7143
$ws.Cells.Item(19, 11).Value = 497.7143
$ws.Cells.Item(19, 13).Value = -322.7143

$ws.Cells.Item(33, 8).Value = 138.91667
$ws.Cells.Item(33, 9).Value = 152.7
$ws.Cells.Item(33, 11).Value = 152.7
$ws.Cells.Item(33, 13).Value = 76.30000000000001

$ws.Cells.Item(40, 8).Value = 6483.3823
$ws.Cells.Item(40, 9).Value = 1683.3334
$ws.Cells.Item(40, 10).Value = 8211.4
$ws.Cells.Item(40, 11).Value = 1683.3334
$ws.Cells.Item(40, 12).Value = 8211.4
$ws.Cells.Item(40, 13).Value = -1508.3334
$ws.Cells.Item(40, 14).Value = -8561.4

$ws.Cells.Item(58, 8).Value = 2203.75
$ws.Cells.Item(58, 9).Value = 1407.5
$ws.Cells.Item(58, 10).Value = 3000
$ws.Cells.Item(58, 11).Value = 4222.5
$ws.Cells.Item(58, 12).Value = 9000
$ws.Cells.Item(58, 13).Value = -4072.5
$ws.Cells.Item(58, 14).Value = -9300

$ws.Cells.Item(100, 8).Value = 3499.6667
$ws.Cells.Item(100, 9).Value = 0
$ws.Cells.Item(100, 10).Value = 3499.6667
$ws.Cells.Item(100, 11).Value = 0
$ws.Cells.Item(100, 12).Value = 3499.6667
$ws.Cells.Item(100, 13).ClearContents()
$ws.Cells.Item(100, 14).Value = -4581.6667

$ws.Cells.Item(113, 8).Value = 4498.6665
$ws.Cells.Item(113, 9).Value = 4498.6665
$ws.Cells.Item(113, 10).Value = 0
$ws.Cells.Item(113, 11).Value = 4498.6665
$ws.Cells.Item(113, 12).Value = 0
$ws.Cells.Item(113, 13).Value = -1244.6665
$ws.Cells.Item(113, 14).ClearContents()

$ws.Cells.Item(129, 8).Value = 842.6
$ws.Cells.Item(129, 9).Value = 842.6
$ws.Cells.Item(129, 10).Value = 0
$ws.Cells.Item(129, 11).Value = 2527.8
$ws.Cells.Item(129, 12).Value = 0
$ws.Cells.Item(129, 13).Value = 2472.2
$ws.Cells.Item(129, 14).ClearContents()

$ws.Cells.Item(132, 8).Value = 1796.5454
$ws.Cells.Item(132, 9).Value = 884.25
$ws.Cells.Item(132, 11).Value = 2652.75
$ws.Cells.Item(132, 13).Value = -122.75

$ws.Cells.Item(135, 8).Value = 1656.4615
$ws.Cells.Item(135, 9).Value = 1558.7273
$ws.Cells.Item(135, 11).Value = 14028.5457
$ws.Cells.Item(135, 13).Value = -11493.5457

$ws.Cells.Item(138, 8).Value = 3512.7144
$ws.Cells.Item(138, 10).Value = 3957.4167
$ws.Cells.Item(138, 12).Value = 11872.2501
$ws.Cells.Item(138, 14).Value = -22152.2501

$ws.Cells.Item(141, 8).Value = 5525.8423
$ws.Cells.Item(141, 9).Value = 5525.8423
$ws.Cells.Item(141, 10).Value = 0
$ws.Cells.Item(141, 11).Value = 16577.5269
$ws.Cells.Item(141, 12).Value = 0
$ws.Cells.Item(141, 13).Value = -11397.5269
$ws.Cells.Item(141, 14).ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(3, 8).Value = 200
$ws.Cells.Item(3, 9).Value = 0
$ws.Cells.Item(3, 10).Value = 200
$ws.Cells.Item(3, 11).Value = 0
$ws.Cells.Item(3, 12).Value = 200
$ws.Cells.Item(3, 13).ClearContents()
$ws.Cells.Item(3, 14).Value = -430

$ws.Cells.Item(132, 8).Value = 1960.375
$ws.Cells.Item(132, 9).Value = 1811.9286
$ws.Cells.Item(132, 11).Value = 5435.7858
$ws.Cells.Item(132, 13).Value = -2905.7858

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(8, 8).Value = 0
$ws.Cells.Item(8, 9).Value = 0
$ws.Cells.Item(8, 10).Value = 0
$ws.Cells.Item(8, 11).Value = 0
$ws.Cells.Item(8, 12).Value = 0
$ws.Cells.Item(8, 13).ClearContents()
$ws.Cells.Item(8, 14).ClearContents()

$ws.Cells.Item(20, 8).Value = 4006
$ws.Cells.Item(20, 9).Value = 4006
$ws.Cells.Item(20, 10).Value = 0
$ws.Cells.Item(20, 11).Value = 4006
$ws.Cells.Item(20, 12).Value = 0
$ws.Cells.Item(20, 13).Value = -3759
$ws.Cells.Item(20, 14).ClearContents()

$ws.Cells.Item(80, 8).Value = 676.2222
$ws.Cells.Item(80, 10).Value = 859.3333
$ws.Cells.Item(80, 12).Value = 859.3333
$ws.Cells.Item(80, 14).Value = -2855.3333

$ws.Cells.Item(82, 8).Value = 7580.1665
$ws.Cells.Item(82, 9).Value = 7580.1665
$ws.Cells.Item(82, 11).Value = 7580.1665
$ws.Cells.Item(82, 13).Value = -7197.1665

$ws.Cells.Item(83, 8).Value = 676.2222
$ws.Cells.Item(83, 10).Value = 859.3333
$ws.Cells.Item(83, 12).Value = 4296.6665
$ws.Cells.Item(83, 14).Value = -14280.6665

$ws.Cells.Item(85, 8).Value = 7580.1665
$ws.Cells.Item(85, 9).Value = 7580.1665
$ws.Cells.Item(85, 11).Value = 7580.1665
$ws.Cells.Item(85, 13).Value = -6254.1665

$ws.Cells.Item(86, 8).Value = 1499.6666
$ws.Cells.Item(86, 9).Value = 1499.6666
$ws.Cells.Item(86, 11).Value = 1499.6666
$ws.Cells.Item(86, 13).Value = -376.6666

$ws.Cells.Item(89, 8).Value = 1499.6666
$ws.Cells.Item(89, 9).Value = 1499.6666
$ws.Cells.Item(89, 11).Value = 7498.333000000001
$ws.Cells.Item(89, 13).Value = -1882.333000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(25, 8).Value = 6837
$ws.Cells.Item(25, 9).Value = 255.5
$ws.Cells.Item(25, 11).Value = 255.5
$ws.Cells.Item(25, 13).Value = -81.5

$ws.Cells.Item(31, 8).Value = 8098
$ws.Cells.Item(31, 9).Value = 10599.75
$ws.Cells.Item(31, 10).Value = 6096.6
$ws.Cells.Item(31, 11).Value = 10599.75
$ws.Cells.Item(31, 12).Value = 6096.6
$ws.Cells.Item(31, 13).Value = -10304.75
$ws.Cells.Item(31, 14).Value = -6686.6

$ws.Cells.Item(34, 8).Value = 8098
$ws.Cells.Item(34, 9).Value = 10599.75
$ws.Cells.Item(34, 10).Value = 6096.6
$ws.Cells.Item(34, 11).Value = 10599.75
$ws.Cells.Item(34, 12).Value = 6096.6
$ws.Cells.Item(34, 13).Value = -10397.75
$ws.Cells.Item(34, 14).Value = -6500.6

$ws.Cells.Item(69, 8).Value = 0
$ws.Cells.Item(69, 9).Value = 0
$ws.Cells.Item(69, 11).Value = 0
$ws.Cells.Item(69, 13).ClearContents()

$ws.Cells.Item(72, 8).Value = 0
$ws.Cells.Item(72, 9).Value = 0
$ws.Cells.Item(72, 11).Value = 0
$ws.Cells.Item(72, 13).ClearContents()

$ws.Cells.Item(105, 8).Value = 1184
$ws.Cells.Item(105, 9).Value = 1210.2858
$ws.Cells.Item(105, 11).Value = 1210.2858
$ws.Cells.Item(105, 13).Value = 536.7141999999999

$ws.Cells.Item(132, 8).Value = 877.5454999999999
$ws.Cells.Item(132, 9).Value = 765.3
$ws.Cells.Item(132, 11).Value = 2295.9
$ws.Cells.Item(132, 13).Value = 234.1000000000004

$ws.Cells.Item(134, 8).Value = 4222.75
$ws.Cells.Item(134, 9).Value = 3297
$ws.Cells.Item(134, 11).Value = 9891
$ws.Cells.Item(134, 13).Value = -7356

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 2000333
$ws.Cells.Item(4, 9).Value = 2500500
$ws.Cells.Item(4, 11).Value = 7501500
$ws.Cells.Item(4, 13).Value = -7501388

$ws.Cells.Item(9, 8).Value = 612.5
$ws.Cells.Item(9, 9).Value = 150
$ws.Cells.Item(9, 10).Value = 2000
$ws.Cells.Item(9, 11).Value = 450
$ws.Cells.Item(9, 12).Value = 6000
$ws.Cells.Item(9, 13).Value = -226
$ws.Cells.Item(9, 14).Value = -6448

$ws.Cells.Item(22, 8).Value = 2000
$ws.Cells.Item(22, 9).Value = 0
$ws.Cells.Item(22, 10).Value = 2000
$ws.Cells.Item(22, 11).Value = 0
$ws.Cells.Item(22, 12).Value = 6000
$ws.Cells.Item(22, 13).ClearContents()
$ws.Cells.Item(22, 14).Value = -6338

$ws.Cells.Item(27, 8).Value = 2000
$ws.Cells.Item(27, 9).Value = 0
$ws.Cells.Item(27, 10).Value = 2000
$ws.Cells.Item(27, 11).Value = 0
$ws.Cells.Item(27, 12).Value = 6000
$ws.Cells.Item(27, 13).ClearContents()
$ws.Cells.Item(27, 14).Value = -6204

$ws.Cells.Item(60, 8).Value = 0
$ws.Cells.Item(60, 10).Value = 0
$ws.Cells.Item(60, 12).Value = 0
$ws.Cells.Item(60, 14).ClearContents()

$ws.Cells.Item(128, 8).Value = 125000
$ws.Cells.Item(128, 9).Value = 125000
$ws.Cells.Item(128, 11).Value = 375000
$ws.Cells.Item(128, 13).Value = -370020

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(113, 8).Value = 0
$ws.Cells.Item(113, 9).Value = 0
$ws.Cells.Item(113, 10).Value = 0
$ws.Cells.Item(113, 11).Value = 0
$ws.Cells.Item(113, 12).Value = 0
$ws.Cells.Item(113, 13).ClearContents()
$ws.Cells.Item(113, 14).ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 1666.3334
$ws.Cells.Item(22, 9).Value = 1499.5
$ws.Cells.Item(22, 10).Value = 2000
$ws.Cells.Item(22, 11).Value = 1499.5
$ws.Cells.Item(22, 12).Value = 2000
$ws.Cells.Item(22, 13).Value = -1204.5
$ws.Cells.Item(22, 14).Value = -2590

$ws.Cells.Item(27, 8).Value = 1666.3334
$ws.Cells.Item(27, 9).Value = 1499.5
$ws.Cells.Item(27, 10).Value = 2000
$ws.Cells.Item(27, 11).Value = 1499.5
$ws.Cells.Item(27, 12).Value = 2000
$ws.Cells.Item(27, 13).Value = -1392.5
$ws.Cells.Item(27, 14).Value = -2214

$ws.Cells.Item(97, 8).Value = 4796.75
$ws.Cells.Item(97, 10).Value = 4796.75
$ws.Cells.Item(97, 12).Value = 4796.75
$ws.Cells.Item(97, 14).Value = -6778.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(44, 8).Value = 12550
$ws.Cells.Item(44, 9).Value = 100
$ws.Cells.Item(44, 10).Value = 25000
$ws.Cells.Item(44, 11).Value = 100
$ws.Cells.Item(44, 12).Value = 25000
$ws.Cells.Item(44, 13).Value = 454
$ws.Cells.Item(44, 14).Value = -26108

$ws.Cells.Item(94, 8).Value = 25999
$ws.Cells.Item(94, 10).Value = 25999
$ws.Cells.Item(94, 12).Value = 25999
$ws.Cells.Item(94, 14).Value = -27801

$ws.Cells.Item(104, 8).Value = 8495
$ws.Cells.Item(104, 10).Value = 8495
$ws.Cells.Item(104, 12).Value = 8495
$ws.Cells.Item(104, 14).Value = -15483

$ws.Cells.Item(132, 8).Value = 3232.7222
$ws.Cells.Item(132, 9).Value = 1560.9231
$ws.Cells.Item(132, 10).Value = 7579.4
$ws.Cells.Item(132, 11).Value = 4682.7693
$ws.Cells.Item(132, 12).Value = 22738.2
$ws.Cells.Item(132, 13).Value = -2152.7693
$ws.Cells.Item(132, 14).Value = -27798.2

$ws.Cells.Item(136, 8).Value = 887.25
$ws.Cells.Item(136, 9).Value = 799.7143
$ws.Cells.Item(136, 10).Value = 1500
$ws.Cells.Item(136, 11).Value = 2399.1429
$ws.Cells.Item(136, 12).Value = 4500
$ws.Cells.Item(136, 13).Value = 150.8571000000002
$ws.Cells.Item(136, 14).Value = -9600
